$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "particip" (E) and "taxa_sucesso" (F) columns were stored as raw
# fractions (e.g. 0.988...) but should be stored as percentage values
# (e.g. 98.8...) -- i.e. each value multiplied by 100.

$ws.Range("E2").Value = 98.80149812734082
$ws.Range("F2").Value = 62.0166793025019

$ws.Range("E3").Value = 1.198501872659176
$ws.Range("F3").Value = 75

$ws.Range("E4").Value = 99.04632152588556
$ws.Range("F4").Value = 94.15405777166438

$ws.Range("E5").Value = 0.9536784741144414
$ws.Range("F5").Value = 100

$ws.Range("E6").Value = 99.85380116959064
$ws.Range("F6").Value = 22.25475841874085

$ws.Range("E7").Value = 0.1461988304093567
